# Update the richness-selection table with revised AICc / delta / logLik values.
# Each "old" value is unique within the document, so a MatchWholeWord
# Find/Replace (one hit at a time, in this exact order) safely targets the
# correct cell even where a later replacement's "new" text happens to equal an
# earlier "old" text (e.g. -6.17, -4.86): those earlier originals are already
# consumed by the time the colliding value is (re)written further down the list.
$d = $word.ActiveDocument

$olds = @(
    "-4.93"
    "21.33"
    "-4.49"
    "23.09"
    "1.75"
    "-6.17"
    "23.81"
    "-6.38"
    "24.23"
    "-5.29"
    "24.67"
    "-5.66"
    "25.41"
    "-4.33"
    "25.53"
    "4.20"
    "-6.15"
    "26.39"
    "-4.86"
    "26.59"
    "-5.82"
    "28.51"
    "-4.45"
    "28.68"
    "-4.67"
    "29.13"
    "-5.06"
    "29.92"
    "8.58"
    "-4.13"
    "31.12"
    "-6.40"
    "32.58"
    "-5.52"
    "33.91"
    "-5.98"
    "34.83"
    "13.50"
    "-4.69"
    "35.48"
)
$news = @(
    "-5.12"
    "21.71"
    "-4.68"
    "23.47"
    "1.76"
    "-6.36"
    "24.18"
    "-6.57"
    "24.60"
    "-5.48"
    "25.05"
    "-5.84"
    "25.79"
    "-4.52"
    "25.92"
    "4.21"
    "-6.33"
    "26.77"
    "-5.05"
    "26.97"
    "-6.01"
    "28.88"
    "-4.63"
    "29.06"
    "-4.86"
    "29.51"
    "-5.26"
    "30.30"
    "8.59"
    "-4.32"
    "31.50"
    "-6.58"
    "32.96"
    "-5.71"
    "34.29"
    "-6.17"
    "35.20"
    "13.49"
    "-4.87"
    "35.86"
)

$notFound = @()
for ($i = 0; $i -lt $olds.Count; $i++) {
    $old = $olds[$i]
    $new = $news[$i]
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                                   $true, 1, $false, $new, 2)
    if (-not $ok) {
        $notFound += $old
    }
}

if ($notFound.Count -gt 0) {
    Write-Output ("NOT FOUND: " + ($notFound -join ", "))
} else {
    Write-Output ("All " + $olds.Count + " replacements applied successfully.")
}
